# Update row 6 (ano=2025) metrics on Sheet1:
#   total_customers (C6):    348 -> 350
#   new_customers   (E6):     67 -> 69
#   new_rate        (G6): 19.25287356321839 -> 19.71428571428572
#   returning_rate  (H6): 80.74712643678161 -> 80.28571428571428
# returning_customers (D6) stays 281; F6 stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C6").Value = 350
$ws.Range("E6").Value = 69
$ws.Range("G6").Value = 19.71428571428572
$ws.Range("H6").Value = 80.28571428571428
